# Update gh-pages to output generated at 456a3b4
# This updates the "F" column (报名/浏览 count, etc.) numeric values across
# the four sheets of the 杭州-漫展信息 workbook: 展览, 演出, 本地生活, 全部类型.
# The "全部类型" sheet is an aggregate of the other three, so matching rows
# there are updated in parallel.

$wb = $excel.ActiveWorkbook

# ---- 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1889
$ws.Range("F8").Value = 259
$ws.Range("F9").Value = 176
$ws.Range("F10").Value = 1125
$ws.Range("F13").Value = 81
$ws.Range("F14").Value = 115
$ws.Range("F16").Value = 244
$ws.Range("F17").Value = 137
$ws.Range("F18").Value = 115
$ws.Range("F19").Value = 1267
$ws.Range("F20").Value = 484
$ws.Range("F21").Value = 172
$ws.Range("F22").Value = 284
$ws.Range("F24").Value = 622
$ws.Range("F25").Value = 1046
$ws.Range("F26").Value = 58
$ws.Range("F28").Value = 2427
$ws.Range("F29").Value = 1193
$ws.Range("F30").Value = 49
$ws.Range("F31").Value = 130
$ws.Range("F32").Value = 338
$ws.Range("F33").Value = 594
$ws.Range("F34").Value = 739
$ws.Range("F35").Value = 812
$ws.Range("F36").Value = 103
$ws.Range("F38").Value = 734
$ws.Range("F39").Value = 240
$ws.Range("F40").Value = 557
$ws.Range("F41").Value = 653
$ws.Range("F42").Value = 288

# ---- 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 171
$ws.Range("F15").Value = 271
$ws.Range("F20").Value = 7
$ws.Range("F22").Value = 15

# ---- 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 858

# ---- 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 858
$ws.Range("F3").Value = 1889
$ws.Range("F10").Value = 259
$ws.Range("F11").Value = 176
$ws.Range("F13").Value = 171
$ws.Range("F14").Value = 1125
$ws.Range("F17").Value = 81
$ws.Range("F18").Value = 115
$ws.Range("F19").Value = 244
$ws.Range("F21").Value = 137
$ws.Range("F22").Value = 116
$ws.Range("F23").Value = 1267
$ws.Range("F24").Value = 484
$ws.Range("F25").Value = 172
$ws.Range("F26").Value = 284
$ws.Range("F28").Value = 1046
$ws.Range("F29").Value = 2427
$ws.Range("F31").Value = 1193
$ws.Range("F32").Value = 49
$ws.Range("F35").Value = 130
$ws.Range("F36").Value = 338
$ws.Range("F37").Value = 594
$ws.Range("F38").Value = 7
$ws.Range("F40").Value = 739
$ws.Range("F41").Value = 812
$ws.Range("F42").Value = 734
$ws.Range("F43").Value = 240
$ws.Range("F44").Value = 557
$ws.Range("F45").Value = 653
$ws.Range("F46").Value = 288
$ws.Range("F47").Value = 15
